# Add two new columns, I (I0) and J (IF), to the sheet, filling in
# header labels and the per-row numeric values, and updating the used
# range (dimension) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Give the new header cells the same style as the existing header cells
# (bold, centered, bordered) by copying the style from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row numeric values for columns I (I0) and J (IF)
$iValues = @(8,1,7,8,1,1,7,6,7,5,6,6,7,6,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$jValues = @(8,4,9,8,4,6,8,9,8,7,8,8,8,8,5,4,5,5,4,5,6,5,6,5,7,7,4,5,3,2,1)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    [void]($ws.Cells.Item($row, 9).Value = $iValues[$idx])
    [void]($ws.Cells.Item($row, 10).Value = $jValues[$idx])
}
